# ============================================================================
# Edit: add "version" (schema version) and "description" columns to the
# "Export as TSV" sheet, with a new "version list" lookup sheet.
#
# Summary of the change (see commit "assays version 1 (#635)"):
#   * A new worksheet "version list" is inserted right after "Export as TSV"
#     containing the single allowed schema version value "1".
#   * Two new leading columns are inserted into "Export as TSV":
#       A = "version"      (validated against 'version list'!$A$1:$A$1)
#       B = "description"  (free text)
#     All the old columns (old A .. old AL) shift two columns to the right
#     (new C .. new AN).
#   * Cell comments on row 1 move along with their columns, and two new
#     comments are added for the new A1/B1 header cells.
# ============================================================================

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Export as TSV")

# ----------------------------------------------------------------------
# Helpers: column-letter <-> column-number conversion.
# ----------------------------------------------------------------------
function ColLettersToNumber([string]$colLetters) {
    $n = 0
    foreach ($ch in $colLetters.ToCharArray()) {
        $n = $n * 26 + ([int][char]$ch - [int][char]'A' + 1)
    }
    return $n
}

function ColNumberToLetters([int]$num) {
    $letters = ""
    $n = $num
    while ($n -gt 0) {
        $rem = ($n - 1) % 26
        $letters = [char](65 + $rem) + $letters
        $n = [int](($n - 1) / 26)
    }
    return $letters
}

# ----------------------------------------------------------------------
# 1. Capture the existing row-1 comments (ORIGINAL ref + text) before we
#    disturb anything. These refs are relative to the CURRENT (pre-insert)
#    layout: A1 .. AL1.
# ----------------------------------------------------------------------
$oldComments = @(
    @{ Ref = 'A1'; Text = 'HuBMAP Display ID of the donor of the assayed tissue.' },
    @{ Ref = 'B1'; Text = 'HuBMAP Display ID of the assayed tissue.' },
    @{ Ref = 'C1'; Text = 'Start date and time of assay, typically a date-time stamped folder generated by the acquisition instrument. YYYY-MM-DD hh:mm, where YYYY is the year, MM is the month with leading 0s, and DD is the day with leading 0s, hh is the hour with leading zeros, mm are the minutes with leading zeros.' },
    @{ Ref = 'D1'; Text = 'DOI for protocols.io referring to the protocol for this assay.' },
    @{ Ref = 'E1'; Text = 'Name of the person responsible for executing the assay.' },
    @{ Ref = 'F1'; Text = 'Email address for the operator.' },
    @{ Ref = 'G1'; Text = 'Name of the principal investigator responsible for the data.' },
    @{ Ref = 'H1'; Text = 'Email address for the principal investigator.' },
    @{ Ref = 'I1'; Text = 'Each assay is placed into one of the following 3 general categories: generation of images of microscopic entities, identification & quantitation of molecules by mass spectrometry, and determination of nucleotide sequence.' },
    @{ Ref = 'J1'; Text = 'The specific type of assay being executed.' },
    @{ Ref = 'K1'; Text = 'Analytes are the target molecules being measured with the assay.' },
    @{ Ref = 'L1'; Text = 'Specifies whether or not a specific molecule(s) is/are targeted for detection/measurement by the assay. The CODEX analyte is protein.' },
    @{ Ref = 'M1'; Text = 'An acquisition instrument is the device that contains the signal detection hardware and signal processing software. Assays generate signals such as light of various intensities or color or signals representing the molecular mass.' },
    @{ Ref = 'N1'; Text = 'Manufacturers of an acquisition instrument may offer various versions (models) of that instrument with different features or sensitivities. Differences in features or sensitivities may be relevant to processing or interpretation of the data.' },
    @{ Ref = 'O1'; Text = 'The manufacturer of the instrument used to prepare the sample for the assay.' },
    @{ Ref = 'P1'; Text = 'The model number/name of the instrument used to prepare the sample for the assay' },
    @{ Ref = 'Q1'; Text = 'DOI for protocols.io referring to the protocol for preparing tissue sections for the assay.' },
    @{ Ref = 'R1'; Text = 'DOI for protocols.io referring to the protocol for preparing reagents for the assay.' },
    @{ Ref = 'S1'; Text = 'Number of mass channels measured' },
    @{ Ref = 'T1'; Text = 'Number of sections' },
    @{ Ref = 'U1'; Text = 'x resolution. Distance between laser ablation shots in the X-dimension.' },
    @{ Ref = 'V1'; Text = 'Units of x resolution distance between laser ablation shots.' },
    @{ Ref = 'W1'; Text = 'y resolution. Distance between laser ablation shots in the Y-dimension.' },
    @{ Ref = 'X1'; Text = 'Units of y resolution distance between laser ablation shots.' },
    @{ Ref = 'Y1'; Text = 'Frequency value of laser ablation (in Hz)' },
    @{ Ref = 'Z1'; Text = 'Frequency unit of laser ablation' },
    @{ Ref = 'AA1'; Text = 'A description of the region of interest (ROI) captured in the image.' },
    @{ Ref = 'AB1'; Text = 'Multiple images (1-n) are acquired from regions of interest (ROI1, ROI2, ROI3, etc) on a slide. The roi_id is a number from 1-n representing the ROI captured on a slide.' },
    @{ Ref = 'AC1'; Text = 'The acquisition_id refers to the directory containing the ROI images for a slide. Together, the acquisition_id and the roi_id indicate the slide-ROI represented in the image.' },
    @{ Ref = 'AD1'; Text = 'Image width value of the ROI acquisition' },
    @{ Ref = 'AE1'; Text = 'Units of image width of the ROI acquisition' },
    @{ Ref = 'AF1'; Text = 'Image height value of the ROI acquisition' },
    @{ Ref = 'AG1'; Text = 'Units of image height of the ROI acquisition' },
    @{ Ref = 'AH1'; Text = 'This refers to the data type, which is a "float" for the IMC counts.' },
    @{ Ref = 'AI1'; Text = 'Type of signal measured per channel (usually dual counts)' },
    @{ Ref = 'AJ1'; Text = 'Relative path to file with antibody information for this dataset.' },
    @{ Ref = 'AK1'; Text = 'Relative path to file with ORCID IDs for contributors for this dataset.' },
    @{ Ref = 'AL1'; Text = 'Relative path to file or directory with instrument data. Downstream processing will depend on filename extension conventions.' }
)

# Remove the old comments outright; we'll re-create them two columns to
# the right after the column insert below.
foreach ($item in $oldComments) {
    $cell = $ws.Range($item.Ref)
    if ($cell.Comment -ne $null) {
        $cell.Comment.Delete() | Out-Null
    }
}

# ----------------------------------------------------------------------
# 2. Insert the new "version list" sheet right after "Export as TSV".
# ----------------------------------------------------------------------
$versionList = $wb.Worksheets.Add($null, $ws)
$versionList.Name = "version list"
$versionList.Range("A1").Value = "1"

# ----------------------------------------------------------------------
# 3. Insert two new columns at the front of "Export as TSV" for
#    "version" and "description". This shifts all existing columns
#    (and their data / data validations) two places to the right
#    automatically; comments must be re-created manually (step 4).
# ----------------------------------------------------------------------
$ws.Range("A1:B1").EntireColumn.Insert() | Out-Null

$ws.Range("A1").Value = "version"
$ws.Range("B1").Value = "description"

# ----------------------------------------------------------------------
# 4. Re-create the row-1 comments, shifted two columns to the right of
#    where they used to be.
# ----------------------------------------------------------------------
foreach ($item in $oldComments) {
    if ($item.Ref -match '^([A-Z]+)(\d+)$') {
        $colLetters = $matches[1]
        $row = $matches[2]
    }
    $newColNum = (ColLettersToNumber $colLetters) + 2
    $newRef = "$(ColNumberToLetters $newColNum)$row"
    $ws.Range($newRef).AddComment($item.Text) | Out-Null
}

# New comments for the two new header cells.
$ws.Range("A1").AddComment('Version of the schema to use when validating this metadata.') | Out-Null
$ws.Range("B1").AddComment('Free-text description of this assay.') | Out-Null

# ----------------------------------------------------------------------
# 5. Add data validation on the new "version" column referencing the new
#    "version list" sheet (matches the pattern used by the other
#    controlled-vocabulary columns in this workbook).
# ----------------------------------------------------------------------
$verRange = $ws.Range("A2:A1048576")
$verRange.Validation.Add(3, 1, 1, "'version list'!`$A`$1:`$A`$1") | Out-Null
$verRange.Validation.ErrorTitle = "Value must come from list"
$verRange.Validation.ErrorMessage = "Value must be one of: 1."
$verRange.Validation.ShowInput = $true
$verRange.Validation.ShowError = $true

Write-Output "Done. Worksheets:"
foreach ($sheet in $wb.Worksheets) {
    Write-Output "  $($sheet.Name)"
}
Write-Output "Header row:"
for ($c = 1; $c -le 40; $c++) {
    $val = $ws.Cells.Item(1, $c).Value2
    if ($val -eq $null) { break }
    Write-Output "  col $c -> $val"
}
